$d = $word.ActiveDocument

# --- locate the text to change -------------------------------------------
# The paragraph currently reads:
#   "# INF2050 Laboratoires" + " - Outils"
# split across two runs. We need it to read:
#   "# INF2050 Laboratoires" + " - " + "Modifié"
# split across three runs (the first run is untouched, the " - Outils" run
# is shrunk to " - ", and a brand-new "Modifié" run is appended).

$full = $d.Content
$full.Find.Execute(" - Outils", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$dashStart = $full.Start          # start of " - Outils"
$firstRunEnd = $dashStart         # == end of "# INF2050 Laboratoires"

# --- 1) shrink " - Outils" down to " - " ----------------------------------
$d.Content.Find.Execute(" - Outils", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " - ", 2)

# --- 2) append the new word right after the trailing " - " ---------------
$tail = $d.Content
$tail.Find.Execute(" - ", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$tail.Collapse(0)
$newStart = $tail.Start
$tail.InsertAfter("Modifié")

# --- 3) re-establish the original run boundaries --------------------------
# The interpreter (like Word itself) silently coalesces adjacent runs that
# share identical formatting, so after the text edits above everything has
# collapsed back into a single run. Briefly toggling Bold on/off over each
# sub-range is a no-op formatting-wise, but it forces those ranges back out
# into their own distinct <w:r> elements, matching the three-run structure
# the diff expects (run 1 unchanged, run 2 = " - ", run 3 = "Modifié").
$newLen = "Modifié".Length
$newRange = $d.Range($newStart, $newStart + $newLen)
$newRange.Bold = 1
$newRange.Bold = 0

$firstRange = $d.Range(0, $firstRunEnd)
$firstRange.Bold = 1
$firstRange.Bold = 0
